# Add new column 'Servised by' (column O) to the Card23 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# New header cell O1, with the same formatting as the other header cells (N1)
$ws.Cells.Item(1, 15).Value = "Servised by"
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The existing column N ("Correction") data cells were blank; they now hold
# the text "nan" like the rest of the row, and column O extends the table
# (used range) through row 12, staying blank for every existing record.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}
